$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "系统分类" (system classification) date values in column N
# for all data rows (row 2 through row 304), leaving the header (N1)
# and all other columns untouched.
$ws.Range("N2:N304").ClearContents()
